# New crime data collected — weekly CompStat (24th Precinct) update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/issue number and the reporting week's date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/19/2024  Through  2/25/2024"

# ---------------------------------------------------------------------
# Helper: write a plain numeric value into a cell, forcing it onto the
# "numeric" style (s=15 in the original workbook, numFmtId #,##0) even
# if the cell previously held a text placeholder value.
# ---------------------------------------------------------------------
function Set-Num($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "#,##0"
    $r.Value = $val
}

# ---------------------------------------------------------------------
# Helper: write a percent-style numeric value (s=16/19/22, numFmtId
# #,##0.0 / #,##0.00) — just a plain value set, style never toggles
# for these columns in this diff.
# ---------------------------------------------------------------------
function Set-Val($addr, $val) {
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------
# Helper: write a text placeholder ("0" or "***.*") into a cell while
# keeping the destination in the workbook's "General" text style
# (s=14) used throughout the sheet for N/A-style entries. We do this by
# borrowing the format from a cell that already has that exact style,
# writing the text through a formula + paste-values round-trip (so the
# string lands as a literal, not auto-coerced back to a number), and
# finally pasting the donor's number format on top.
# ---------------------------------------------------------------------
function Set-Placeholder($addr, $text, $formatDonor) {
    $ws.Range($formatDonor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Formula = "=""" + $text + """"
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 15 — Rape
# ---------------------------------------------------------------------
Set-Num   "C15" 1
Set-Val   "D15" 1
Set-Val   "E15" 0
Set-Num   "F15" 2
Set-Num   "G15" 2
Set-Val   "H15" 0
Set-Num   "I15" 4
Set-Num   "J15" 3
Set-Val   "K15" 33.333333333333
Set-Val   "L15" 300
Set-Val   "M15" 300
Set-Val   "N15" 0

# ---------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------
Set-Num "C16" 5
Set-Num "D16" 1
Set-Val "E16" 400
Set-Num "F16" 18
Set-Num "G16" 7
Set-Val "H16" 157.142857142857
Set-Num "I16" 28
Set-Num "J16" 21
Set-Val "K16" 33.333333333333
Set-Val "L16" 0
Set-Val "M16" -9.677419354838
Set-Val "N16" -73.584905660377

# ---------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------
Set-Num "C17" 1
Set-Num "D17" 1
Set-Val "E17" 0
Set-Num "F17" 8
Set-Num "G17" 9
Set-Val "H17" -11.111111111111
Set-Num "I17" 18
Set-Num "J17" 19
Set-Val "K17" -5.263157894736
Set-Val "L17" -25
Set-Val "M17" 0
Set-Val "N17" -67.272727272727

# ---------------------------------------------------------------------
# Row 18 — Burglary
# ---------------------------------------------------------------------
Set-Num "C18" 6
Set-Num "D18" 5
Set-Val "E18" 20
Set-Num "F18" 11
Set-Num "G18" 8
Set-Val "H18" 37.5
Set-Num "I18" 24
Set-Num "J18" 18
Set-Val "K18" 33.333333333333
Set-Val "L18" -31.428571428571
Set-Val "M18" 9.090909090909
Set-Val "N18" -86.516853932584

# ---------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------
Set-Num "C19" 7
Set-Num "D19" 12
Set-Val "E19" -41.666666666666
Set-Num "F19" 33
Set-Num "G19" 38
Set-Val "H19" -13.157894736842
Set-Num "I19" 62
Set-Num "J19" 86
Set-Val "K19" -27.906976744186
Set-Val "L19" -17.333333333333
Set-Val "M19" -15.068493150684
Set-Val "N19" -56.338028169014

# ---------------------------------------------------------------------
# Row 20 — G.L.A.   (C20 toggles number -> text placeholder "0")
# ---------------------------------------------------------------------
Set-Placeholder "C20" "0" "C27"
Set-Num "D20" 2
Set-Val "E20" -100
Set-Num "F20" 5
Set-Num "G20" 9
Set-Val "H20" -44.444444444444
Set-Num "I20" 8
Set-Num "J20" 20
Set-Val "K20" -60
Set-Val "L20" 14.285714285714
Set-Val "M20" 166.666666666667
Set-Val "N20" -94.202898550724

# ---------------------------------------------------------------------
# Row 21 — TOTAL
# ---------------------------------------------------------------------
Set-Num "C21" 20
Set-Num "D21" 22
Set-Val "E21" -9.090909090909
Set-Num "F21" 77
Set-Num "G21" 74
Set-Val "H21" 4.054054054054
Set-Num "I21" 144
Set-Num "J21" 169
Set-Val "K21" -14.792899408284
Set-Val "L21" -15.294117647058
Set-Val "M21" -3.355704697986
Set-Val "N21" -77.033492822966

# ---------------------------------------------------------------------
# Row 22 — Transit   (C22 toggles number -> text placeholder "0")
# ---------------------------------------------------------------------
Set-Placeholder "C22" "0" "C27"
Set-Num "F22" 3
Set-Num "G22" 1
Set-Val "H22" 200
Set-Val "L22" -14.285714285714

# ---------------------------------------------------------------------
# Row 23 — Housing
# ---------------------------------------------------------------------
Set-Num "C23" 2
Set-Num "D23" 1
Set-Val "E23" 100
Set-Num "G23" 12
Set-Val "H23" -33.333333333333
Set-Num "I23" 14
Set-Num "J23" 16
Set-Val "K23" -12.5
Set-Val "L23" -6.666666666666
Set-Val "M23" -22.222222222222

# ---------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------
Set-Num "C24" 34
Set-Num "D24" 31
Set-Val "E24" 9.677419354838
Set-Num "F24" 101
Set-Num "G24" 134
Set-Val "H24" -24.626865671641
Set-Num "I24" 211
Set-Num "J24" 326
Set-Val "K24" -35.276073619631
Set-Val "L24" -25.441696113074
Set-Val "M24" 40.666666666666

# ---------------------------------------------------------------------
# Row 25 — Misd. Assault
# ---------------------------------------------------------------------
Set-Num "C25" 2
Set-Num "D25" 1
Set-Val "E25" 100
Set-Num "F25" 17
Set-Num "G25" 16
Set-Val "H25" 6.25
Set-Num "I25" 42
Set-Num "J25" 40
Set-Val "K25" 5
Set-Val "L25" -2.325581395348
Set-Val "M25" -8.695652173913

# ---------------------------------------------------------------------
# Row 26 — UCR Rape*   (C26 toggles text placeholder "0" -> number)
# ---------------------------------------------------------------------
Set-Num "C26" 1
Set-Num "D26" 1
Set-Val "E26" 0
Set-Num "F26" 2
Set-Num "G26" 3
Set-Val "H26" -33.333333333333
Set-Num "I26" 4
Set-Num "J26" 4
Set-Val "K26" 0
Set-Val "L26" 100

# ---------------------------------------------------------------------
# Row 27 — Other Sex Crimes
# ---------------------------------------------------------------------
Set-Num "G27" 2
Set-Val "H27" 50
Set-Val "L27" -37.5

# ---------------------------------------------------------------------
# Row 30 — Hate Crimes   (F30 toggles number -> text placeholder "0")
# ---------------------------------------------------------------------
Set-Placeholder "F30" "0" "C27"
Set-Val "L30" -25

$excel.CutCopyMode = 0
